# Update cryptocurrency price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.338.53'
$ws.Range('D3').Value = '2.657.17'
$ws.Range('E3').Value = '  +3.55%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'609.51"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.14%  '
$ws.Range('D6').Value = "'143.74"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -0.42%  '
$ws.Range('D9').Value = '2.655.75'
$ws.Range('E9').Value = '  +3.54%  '
$ws.Range('E10').Value = '  +1.16%  '
$ws.Range('D11').Value = "'5.62"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.43%  '
$ws.Range('E12').Value = '  +0.34%  '
$ws.Range('D13').Value = "'0.361"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.43%  '
$ws.Range('D14').Value = "'27.33"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.14%  '
$ws.Range('D15').Value = '3.132.17'
$ws.Range('E15').Value = '  +3.36%  '
$ws.Range('D16').Value = '63.247.51'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').Value = "'0.0000144"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').Value = '2.665.48'
$ws.Range('E18').Value = '  +3.73%  '
$ws.Range('D19').Value = "'11.44"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.96%  '
$ws.Range('D20').Value = "'341.77"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('E21').Value = '  +2.38%  '
$ws.Range('D22').Value = "'6.86"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.78%  '
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('E24').Value = '  -1.25%  '
$ws.Range('E25').Value = '  +3.03%  '
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('D27').Value = "'8.65"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.07%  '
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').Value = "'546.39"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +17.10%  '
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('D31').Value = "'7.82"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.06%  '
$ws.Range('E32').Value = '  +5.68%  '
$ws.Range('E33').Value = '  +7.49%  '
$ws.Range('D34').Value = '0.0₃0807'
$ws.Range('E34').Value = '  +1.47%  '
$ws.Range('D35').Value = "'173.03"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.67%  '
$ws.Range('E36').Value = '  +13.83%  '
$ws.Range('E37').Value = '  +2.68%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = "'19.13"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.83%  '
$ws.Range('E40').Value = '  +9.77%  '
$ws.Range('D41').Value = "'174.58"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +11.14%  '
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('E43').Value = '  +2.31%  '
$ws.Range('D44').Value = "'22.20"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.46%  '
$ws.Range('E45').Value = '  +6.65%  '
$ws.Range('D46').Value = "'0.632"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.49%  '
$ws.Range('D47').Value = "'0.0961"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.17%  '
$ws.Range('D48').Value = "'0.0239"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('D49').Value = "'18.74"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.87%  '
$ws.Range('D50').Value = "'1.75"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.29%  '
$ws.Range('E51').Value = '  -0.86%  '
